$wb = $excel.ActiveWorkbook

# --- plotConfiguration sheet: add "subtitle" column and example title/subtitle values ---
$ws2 = $wb.Worksheets.Item("plotConfiguration")
$ws2.Range("O1").Value = "subtitle"
$ws2.Range("D2").Value = "PlotTitle"
$ws2.Range("O2").Value = "PlotSubtitle"

# --- plotGrids sheet: add "subtitle" column and example title/subtitle values ---
$ws3 = $wb.Worksheets.Item("plotGrids")
$ws3.Range("D1").Value = "subtitle"
$ws3.Range("C2").Value = "GridTitle"
$ws3.Range("D2").Value = "GridSubtitle"

# --- Selections / active sheet, matching the final workbook view state ---
$ws2.Range("O3").Select()
$ws3.Activate()
$ws3.Range("D3").Select()
